# Adding showPorPeriodoActividad to periodoController
# Applies the content edits to FIN13Inicial.xlsx:
#  - Sheet "Hoja1. Actividades": group/director/semester header values,
#    month-letter header row, a test "Funciona" row, and clearing a
#    leftover test data row.
#  - Sheet "Hoja2. Integrantes": clearing leftover test member rows.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Hoja1. Actividades")
$ws2 = $wb.Worksheets.Item("Hoja2. Integrantes")

# --- Hoja1. Actividades ---------------------------------------------------

# Group name field (was a stray numeric 1, now the test group name)
$ws1.Range("E4").Value = "qqqqqqqqqq"

# Director name
$ws1.Range("E5").Value = "Prueba Coordinador"

# Semester
$ws1.Range("E6").Value = "2021-1"

# CRONOGRAMA month-letter headers (row 9: D..H)
$ws1.Range("D9").Value = "F"
$ws1.Range("E9").Value = "M"
$ws1.Range("F9").Value = "A"
$ws1.Range("G9").Value = "M"
$ws1.Range("H9").Value = "J"

# Row 10 - sample activity row
$ws1.Range("A10").Value = "Funciona"
$ws1.Range("B10").Value = "Funciona"
$ws1.Range("C10").Value = "Funciona"
$ws1.Range("G10").Value = "X"
$ws1.Range("I10").Value = "Funciona"

# Row 11 - clear leftover test activity row
$ws1.Range("A11").Value = ""
$ws1.Range("B11").Value = ""
$ws1.Range("C11").Value = ""
$ws1.Range("G11").Value = ""
$ws1.Range("I11").Value = ""
$ws1.Range("J11").Value = ""

# --- Hoja2. Integrantes ----------------------------------------------------

# Rows 6 and 7 - clear leftover test member rows
$ws2.Range("A6").Value = ""
$ws2.Range("B6").Value = ""
$ws2.Range("C6").Value = ""
$ws2.Range("D6").Value = ""
$ws2.Range("E6").Value = ""

$ws2.Range("A7").Value = ""
$ws2.Range("B7").Value = ""
$ws2.Range("C7").Value = ""
$ws2.Range("D7").Value = ""
$ws2.Range("E7").Value = ""
